$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# --- Row 2: DOI 10.1302/0301-620X.102B9.BJJ-2020-1127.R1 ---
$c2 = @'
"The effects of COVID-19 on perioperative morbidity and mortality in patients with hip fractures"
'@
$d2 = @'
"Aims During the COVID-19 pandemic, many patients continue to require urgent surgery for hip fractures.
 However, the impact of COVID-19 on perioperative outcomes in these high-risk patients remains unknown.
 The objectives of this study were to establish the effects of COVID-19 on perioperative morbidity and mortality, and determine any risk factors for increased mortality in patients with COVID-19 undergoing hip fracture surgery.
 Methods This multicentre cohort study included 340 COVID-19-negative patients versus 82 COVID-19-positive patients undergoing surgical treatment for hip fractures across nine NHS hospitals in Greater London, UK.
 Patients in both treatment groups were comparable for age, sex, body mass index, fracture configuration, and type of surgery performed.
 Predefined perioperative outcomes were recorded within a 30-day postoperative period.
 Univariate and multivariate analysis were used to identify risk factors associated with increased risk of mortality.
 Results COVID-19-positive patients had increased postoperative mortality rates (30.5% (25\/82) vs 10.3% (35\/340) respectively, p &lt; 0.001) compared to COVID-19-negative patients.
 Risk factors for increased mortality in patients with COVID-19 undergoing surgery included positive smoking status (hazard ratio (HR) 15.4 (95% confidence interval (CI) 4.55 to 52.2; p &lt; 0.001) and greater than three comorbidities (HR 13.5 (95% CI 2.82 to 66.0, p &lt; 0.001).
 COVID-19-positive patients had increased risk of postoperative complications (89.0% (73\/82) vs 35.0% (119\/340) respectively; p &lt; 0.001), more critical care unit admissions (61.0% (50\/82) vs 18.2% (62\/340) respectively; p &lt; 0.001), and increased length of hospital stay (mean 13.8 days (SD 4.6) vs 6.7 days (SD 2.5) respectively; p &lt; 0.001), compared to COVID-19-negative patients.
 Conclusion Hip fracture surgery in COVID-19-positive patients was associated with increased length of hospital stay, more admissions to the critical care unit, higher risk of perioperative complications, and increased mortality rates compared to COVID-19-negative patients.
 Risk factors for increased mortality in patients with COVID-19 undergoing surgery included positive smoking status and multiple (greater than three) comorbidities.
 Cite this article: Bone Joint J 2020;102-B(9):1136\u20131145. 
'@
$e2 = @'
[Babar%Kayani%xref no email%0,  Elliot%Onochie%xref no email%0,  Vijay%Patil%xref no email%0,  Fahima%Begum%xref no email%0,  Rory%Cuthbert%xref no email%0,  David%Ferguson%xref no email%0,  Jagmeet S.%Bhamra%xref no email%0,  Aadhar%Sharma%xref no email%0,  Peter%Bates%xref no email%0,  Fares S.%Haddad%xref no email%0]
'@
$ws.Range("C2").Value = $c2
$ws.Range("D2").Value = $d2
$ws.Range("E2").Value = $e2
$ws.Range("F2").Value = $ws.Range("A2").Value2
$ws.Range("G2").Value = "CROSSREF"
$ws.Range("H2").Value = "2023-05-26"

# --- Row 3: Authors whitespace reformat ---
$e3 = @'
[Kenneth A.%Egol%NULL%0,         Sanjit R.%Konda%NULL%0,         Mackenzie L.%Bird%NULL%0,         Nicket%Dedhia%NULL%0,         Emma K.%Landes%NULL%0,         Rachel A.%Ranson%NULL%0,         Sara J.%Solasz%NULL%0,         Vinay K.%Aggarwal%NULL%0,         Joseph A.%Bosco%NULL%0,         David L.%Furgiuele%NULL%0,         Abhishek%Ganta%NULL%0,         Jason%Gould%NULL%0,         Thomas R.%Lyon%NULL%0,         Toni M.%McLaurin%NULL%0,         Nirmal C.%Tejwani%NULL%0,         Joseph D.%Zuckerman%NULL%0,         Philipp%Leucht%NULL%0]
'@
$ws.Range("E3").Value = $e3

# --- Row 4: Authors whitespace reformat ---
$e4 = @'
[Drake G.%LeBrun%NULL%0,         Maxwell A.%Konnaris%NULL%0,         Gregory C.%Ghahramani%NULL%0,         Ajay%Premkumar%NULL%0,         Chris J.%DeFrancesco%NULL%0,         Jordan A.%Gruskay%NULL%0,         Aleksey%Dvorzhinskiy%NULL%0,         Milan S.%Sandhu%NULL%0,         Elan M.%Goldwyn%NULL%0,         Christopher L.%Mendias%NULL%0,         William M.%Ricci%NULL%0]
'@
$ws.Range("E4").Value = $e4

# --- Row 5: Authors whitespace reformat ---
$e5 = @'
[Amit%Thakrar%NULL%1,         Karen%Chui%NULL%1,         Akhil%Kapoor%NULL%1,         John%Hambidge%NULL%1]
'@
$ws.Range("E5").Value = $e5

# --- Row 6: DOI 10.1302/0301-620X.102B9.BJJ-2020-1100.R1 ---
$c6 = @'
"IMPACT-Scot report on COVID-19 and hip fractures"
'@
$d6 = @'
"Aims The primary aim was to assess the independent influence of coronavirus disease (COVID-19) on 30-day mortality for patients with a hip fracture.
 The secondary aims were to determine whether: 1) there were clinical predictors of COVID-19 status; and 2) whether social lockdown influenced the incidence and epidemiology of hip fractures.
 Methods A national multicentre retrospective study was conducted of all patients presenting to six trauma centres or units with a hip fracture over a 46-day period (23 days pre- and 23 days post-lockdown).
 Patient demographics, type of residence, place of injury, presentation blood tests, Nottingham Hip Fracture Score, time to surgery, operation, American Society of Anesthesiologists (ASA) grade, anaesthetic, length of stay, COVID-19 status, and 30-day mortality were recorded.
 Results Of 317 patients with acute hip fracture, 27 (8.5%) had a positive COVID-19 test.
 Only seven (26%) had suggestive symptoms on admission.
 COVID-19-positive patients had a significantly lower 30-day survival compared to those without COVID-19 (64.5%, 95% confidence interval (CI) 45.7 to 83.3 vs 91.7%, 95% CI 88.2 to 94.8; p &lt; 0.001).
 COVID-19 was independently associated with increased 30-day mortality risk adjusting for: 1) age, sex, type of residence (hazard ratio (HR) 2.93; p = 0.008); 2) Nottingham Hip Fracture Score (HR 3.52; p = 0.001); and 3) ASA (HR 3.45; p = 0.004).
 Presentation platelet count predicted subsequent COVID-19 status; a value of &lt; 217 \u00d7 109\/l was associated with 68% area under the curve (95% CI 58 to 77; p = 0.002) and a sensitivity and specificity of 63%.
 A similar number of patients presented with hip fracture in the 23 days pre-lockdown (n = 160) and 23 days post-lockdown (n = 157) with no significant (all p \u2265 0.130) difference in patient demographics, residence, place of injury, Nottingham Hip Fracture Score, time to surgery, ASA, or management.
 Conclusion COVID-19 was independently associated with an increased 30-day mortality rate for patients with a hip fracture.
 Notably, most patients with hip fracture and COVID-19 lacked suggestive symptoms at presentation.
 Platelet count was an indicator of risk of COVID-19 infection.
 These findings have implications for the management of hip fractures, in particular the need for COVID-19 testing.
 Cite this article: Bone Joint J 2020;102-B(9):1219\u20131228. 
'@
$e6 = @'
[Andrew J.%Hall%xref no email%0,  Nicholas D.%Clement%xref no email%0,  Luke%Farrow%xref no email%0,  Alasdair M. J.%MacLullich%xref no email%0,  Graham F.%Dall%xref no email%0,  Chloe E. H.%Scott%xref no email%0,  Paul J.%Jenkins%xref no email%0,  Timothy O.%White%xref no email%0,  Andrew D.%Duckworth%xref no email%0]
'@
$ws.Range("C6").Value = $c6
$ws.Range("D6").Value = $d6
$ws.Range("E6").Value = $e6
$ws.Range("F6").Value = $ws.Range("A6").Value2
$ws.Range("G6").Value = "CROSSREF"
$ws.Range("H6").Value = "2023-05-26"

# --- Row 7: Authors whitespace reformat ---
$e7 = @'
[Karen%Chui%NULL%1,         Amit%Thakrar%NULL%1,         Shivakumar%Shankar%NULL%1]
'@
$ws.Range("E7").Value = $e7
